$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent results for "case with 380 kV" (commit message).
# Columns B, C, D, E, G, K, L, O get new values for data rows 2-25 (index 0-23 in col A).
# Columns F, H, I, J, M, N remain 0 and column A (index) is unchanged.

$colBValues = @(13.94166074358692, 13.74337575544873, 13.62357093761007, 13.5753011075133, 13.56732108111844, 13.62291763767345, 13.87292683553206, 14.3758177311754, 14.74925228347267, 14.91921949018864, 14.98353152837749, 14.96968403696767, 14.92451185795646, 14.89683411095084, 14.73814084107026, 14.64076472884298, 14.58476874027813, 14.56581352603288, 14.65112979017437, 14.93778190084606, 15.12480105217829, 15.02503566867447, 14.64644378205828, 14.23880877677607)
for ($i = 0; $i -lt $colBValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colBValues[$i]
}

$colCValues = @(8.672147381132019, 8.662557816197673, 8.6565733223243, 8.65411062757912, 8.653700260145877, 8.656540206206618, 8.668861091894268, 8.69224505033605, 8.708936809201392, 8.7164212644812, 8.719239521413764, 8.718633273402832, 8.71665344549082, 8.715438658157842, 8.70844546825322, 8.704127307300421, 8.701633401578658, 8.700787263255751, 8.70458804225586, 8.717235404124308, 8.725408110229406, 8.721054802009046, 8.704379779277222, 8.686003795653866)
for ($i = 0; $i -lt $colCValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colCValues[$i]
}

$colDValues = @(6.019309926618026, 5.906564960196289, 5.838037020937293, 5.810328858099723, 5.805742307582011, 5.837662402554254, 5.980315753530063, 6.263768197616939, 6.471889486011007, 6.566058340675488, 6.601606805611062, 6.593956345772502, 6.568985376714791, 6.553674344057553, 6.465721811376671, 6.411608126826318, 6.380437868952324, 6.369877541603071, 6.417373613931508, 6.576323271019305, 6.679541263376668, 6.62452488876484, 6.414767219869745, 6.186941029675165)
for ($i = 0; $i -lt $colDValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $colDValues[$i]
}

$colEValues = @(12.0832549661489, 12.08956133031867, 12.09554984251218, 12.09852261817932, 12.09904840463654, 12.09558777853193, 12.08499033206998, 12.08098539561807, 12.08824028254974, 12.09374469136431, 12.09614477522928, 12.09561385645382, 12.09393582891829, 12.0929490588566, 12.08792480549368, 12.0854062815476, 12.08416524456063, 12.08378073315558, 12.08565291046319, 12.09442014982489, 12.10198936191314, 12.09778170486006, 12.08554076498256, 12.08027512795514)
for ($i = 0; $i -lt $colEValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $colEValues[$i]
}

$colGValues = @(3.69075369690394, 3.693074343023472, 3.69457466832144, 3.695205095064807, 3.695310928261838, 3.694583093335327, 3.691538235784533, 3.686163017443916, 3.682573049910634, 3.68101703330083, 3.680438829744564, 3.680562866744848, 3.680969243493506, 3.681219595217666, 3.682676285328796, 3.68358961790865, 3.684122200972728, 3.684303772829598, 3.683491641304684, 3.68084958200101, 3.679187086050999, 3.680068532121219, 3.683535913173102, 3.687553790783222)
for ($i = 0; $i -lt $colGValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $colGValues[$i]
}

$colKValues = @(10.28108107126131, 10.14738423501728, 10.06672137193081, 10.03424977300274, 10.02888316016758, 10.06628177971173, 10.23471041353696, 10.57453628177276, 10.82761611732511, 10.94298415852409, 10.98666439810569, 10.97725804785174, 10.94657813653089, 10.92778365448632, 10.82007779682357, 10.75403538459766, 10.71607463848067, 10.70322738206123, 10.76106341977744, 10.95559009134754, 11.08266434994473, 11.01486144561945, 10.75788601721065, 10.48183222067704)
for ($i = 0; $i -lt $colKValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 11).Value = $colKValues[$i]
}

$colLValues = @(9.780198257935977, 9.770369215596867, 9.766082296536061, 9.764776653477316, 9.764586551794112, 9.766062899240307, 9.776447484736234, 9.810596509211576, 9.843959408931253, 9.860900454486753, 9.867566001538554, 9.866119378588383, 9.861443835869981, 9.85861242933221, 9.842887490824843, 9.833690288405439, 9.828566562086671, 9.826860409228102, 9.834652162213992, 9.862810389782117, 9.882670868951903, 9.871938768736541, 9.834216788388352, 9.79989487893949)
for ($i = 0; $i -lt $colLValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 12).Value = $colLValues[$i]
}

$colOValues = @(29.97974670423154, 30.01390643766295, 30.04064867267767, 30.05299454173044, 30.05513195788085, 30.04080931301352, 29.99032660703076, 29.93718177317068, 29.92618313212186, 29.92728064853561, 29.92857352136514, 29.92825606723685, 29.9273694341956, 29.9269405824061, 29.92623418081467, 29.92736354243922, 29.92858739614296, 29.92910039681853, 29.92718388322649, 29.92760605378081, 29.93299497699937, 29.92965110780917, 29.9272633173448, 29.94663853101006)
for ($i = 0; $i -lt $colOValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 15).Value = $colOValues[$i]
}
